# Apply the "Updated Master Data excels" edit:
#  - remove the helper "Sheet1" worksheet (and its Table1, which lives on it)
#  - append new template-type master-data rows (92-121) to master-template_type
#  - bump vertical DPI on the page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-template_type")

# --- Delete the helper "Sheet1" tab (this also drops the Table1 defined on it) ---
$excel.DisplayAlerts = $false
$helper = $wb.Worksheets.Item("Sheet1")
$helper.Delete()
$excel.DisplayAlerts = $true

# --- Append the new master-data rows starting at row 92 ---
$rows = @(
    @("RPR_UIN_CARD_TEMPLATE", "UIN card template", "eng", $true, "superadmin", "now()"),
    @("RPR_UIN_CARD_TEMPLATE", "قالب بطاقة UIN", "ara", $true, "superadmin", "now()"),
    @("RPR_UIN_CARD_TEMPLATE", "Modèle de carte UIN", "fra", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_SMS", "Template for UIN Deactivation SMS", "eng", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_SMS", "قالب لتعطيل UIN SMS", "ara", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_SMS", "Modèle pour SMS de désactivation UIN", "fra", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_EMAIL", "Template for UIN Deactivation Email", "eng", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_EMAIL", "قالب لإلغاء تنشيط البريد", "ara", $true, "superadmin", "now()"),
    @("RPR_UIN_DEAC_EMAIL", "Modèle pour Email de désactivation UIN", "fra", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_SMS", "Template for UIN Reactivate SMS", "eng", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_SMS", "قالب لـ UIN تنشيط SMS", "ara", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_SMS", "Modèle pour UIN Réactiver SMS", "fra", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_EMAIL", "Template for UIN Reactivate Email", "eng", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_EMAIL", "قالب لـ UIN تنشيط البريد", "ara", $true, "superadmin", "now()"),
    @("RPR_UIN_REAC_EMAIL", "Modèle pour UIN Réactiver Email", "fra", $true, "superadmin", "now()"),
    @("reg-sms-notification", "Registration Acknowledgement Template", "eng", $true, "superadmin", "now()"),
    @("reg-sms-notification", "نموذج شكر التسجيل", "ara", $true, "superadmin", "now()"),
    @("reg-sms-notification", "accusé de réception", "fra", $true, "superadmin", "now()"),
    @("reg-email-notification", "Registration Acknowledgement Template", "eng", $true, "superadmin", "now()"),
    @("reg-email-notification", "نموذج شكر التسجيل", "ara", $true, "superadmin", "now()"),
    @("reg-email-notification", "accusé de réception", "fra", $true, "superadmin", "now()"),
    @("reg-ack-template-part1", "Registration Acknowledgement Template - Part 1", "eng", $true, "superadmin", "now()"),
    @("reg-ack-template-part2", "نموذج شكر التسجيل", "ara", $true, "superadmin", "now()"),
    @("reg-ack-template-part3", "accusé de réception", "fra", $true, "superadmin", "now()"),
    @("reg-ack-template-part2", "Registration Acknowledgement Template - Part 2", "eng", $true, "superadmin", "now()"),
    @("reg-ack-template-part3", "نموذج شكر التسجيل", "ara", $true, "superadmin", "now()"),
    @("reg-ack-template-part4", "accusé de réception", "fra", $true, "superadmin", "now()"),
    @("reg-ack-template-part3", "Registration Acknowledgement Template - Part 3", "eng", $true, "superadmin", "now()"),
    @("reg-ack-template-part4", "نموذج شكر التسجيل", "ara", $true, "superadmin", "now()"),
    @("reg-ack-template-part5", "accusé de réception", "fra", $true, "superadmin", "now()")
)

$startRow = 92
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# --- Misc formatting bump recorded in the diff ---
$ws.PageSetup.VerticalDpi = 300

# --- Move the active selection to just below the newly-added data, matching
#     the author's final cursor position when they saved the sheet ---
$ws.Range("A122:XFD1048576").Select()
